# Append one new data row (row 43) to each of the four worksheets, mirroring
# the existing table layout (columns A-I). This extends each sheet's used
# range from A1:I42 to A1:I43.

$wb = $excel.ActiveWorkbook

$newRows = @{
    "ROW35-FE-LIFTER" = @{
        A = 45745.83098565972
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x76"
        E = "0xd"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 374
        I = 13
    }
    "ROW35-MID-LIFTER" = @{
        A = 45745.67990320602
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x7a"
        E = "0xe"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 378
        I = 14
    }
    "ROW02-FE-LIFTER" = @{
        A = 45745.82222035879
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x76"
        E = "0x3"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 374
        I = 3
    }
    "ROW02-MID-LIFTER" = @{
        A = 45745.8802064699
        B = "0x01,0x90"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x01,0x7a"
        E = "0x3"
        F = 400
        G = [double]"9.85046333984776e+23"
        H = 378
        I = 3
    }
}

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $row = $newRows[$ws.Name]
    if ($row -eq $null) { continue }

    $targetRow = 43

    $ws.Cells.Item($targetRow, 1).Value = $row.A
    $ws.Cells.Item($targetRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($targetRow, 2).Value = $row.B
    $ws.Cells.Item($targetRow, 3).Value = $row.C
    $ws.Cells.Item($targetRow, 4).Value = $row.D
    $ws.Cells.Item($targetRow, 5).Value = $row.E
    $ws.Cells.Item($targetRow, 6).Value = $row.F
    $ws.Cells.Item($targetRow, 7).Value = $row.G
    $ws.Cells.Item($targetRow, 8).Value = $row.H
    $ws.Cells.Item($targetRow, 9).Value = $row.I
}
